$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.119.04'
$ws.Range("E2").Value = '  -0.12%  '
$ws.Range("D3").Value = '2.321.48'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '303.46'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '99.76'
$ws.Range("D6").ClearFormats()
$ws.Range("E7").Value = '  +0.70%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.519'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +2.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.13'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.88%  '
$ws.Range("E11").Value = '  -0.65%  '
$ws.Range("E12").Value = '  -0.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.77'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.92'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.29%  '
$ws.Range("D15").Value = '2.682.61'
$ws.Range("E15").Value = '  +0.48%  '
$ws.Range("D16").Value = '2.339.00'
$ws.Range("E16").Value = '  +3.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.796'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.89%  '
$ws.Range("D18").Value = '43.034.47'
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.28'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +4.86%  '
$ws.Range("E20").Value = '  +1.22%  '
$ws.Range("D21").Value = '0.0₃0909'
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("E22").Value = '  +0.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '240.04'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +1.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.17'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.46'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.38%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.48'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.76%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '169.01'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.07'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("E31").Value = '  -9.29%  '
$ws.Range("E32").Value = '  +8.17%  '
$ws.Range("E33").Value = '  +2.73%  '
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.11'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +5.55%  '
$ws.Range("E36").Value = '  -1.54%  '
$ws.Range("E37").Value = '  +0.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.84'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.83%  '
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("E40").Value = '  -1.20%  '
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("D42").Value = '1.994.09'
$ws.Range("E43").Value = '  +1.19%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.17'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -8.18%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.22'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.75%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.63'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("E47").Value = '  +0.02%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.17'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.57%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '76.24'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +8.39%  '
$ws.Range("D50").Value = '2.547.55'
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("E51").Value = '  +0.53%  '
